$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.03928371021970611
$ws.Cells.Item(2, 8).Value = -7.083831197969451
$ws.Cells.Item(2, 9).Value = 45.21246975919157
$ws.Cells.Item(3, 7).Value = 0.06120378081178399
$ws.Cells.Item(3, 8).Value = 25.67421344941646
$ws.Cells.Item(4, 7).Value = -0.018539108938213
$ws.Cells.Item(4, 8).Value = -1073.746112098445
$ws.Cells.Item(5, 7).Value = 0.06238096074884684
$ws.Cells.Item(5, 8).Value = 1543.182150897558
$ws.Cells.Item(6, 7).Value = 0.02931899468867842
$ws.Cells.Item(6, 8).Value = -15.43308240821518
$ws.Cells.Item(7, 7).Value = 0.03544059348977226
$ws.Cells.Item(7, 8).Value = -33.36849880243554
$ws.Cells.Item(8, 7).Value = -0.01037161503968795
$ws.Cells.Item(8, 8).Value = 44.8927314306451
$ws.Cells.Item(9, 7).Value = -0.02599485212289188
$ws.Cells.Item(9, 8).Value = -20.74956715828752
$ws.Cells.Item(10, 7).Value = -0.06970307957957617
$ws.Cells.Item(10, 8).Value = 4.124099723280551
$ws.Cells.Item(11, 7).Value = -0.06793821490253291
$ws.Cells.Item(11, 8).Value = 26.16297590730446
$ws.Cells.Item(12, 7).Value = -0.2318843439180897
$ws.Cells.Item(12, 8).Value = 5.14218101998676
$ws.Cells.Item(13, 7).Value = -0.2949552863068522
$ws.Cells.Item(13, 8).Value = -7.329365491908443
$ws.Cells.Item(14, 7).Value = -0.07248458257264322
$ws.Cells.Item(14, 8).Value = -95.39101350386515
$ws.Cells.Item(15, 7).Value = 0.001293213872535024
$ws.Cells.Item(15, 8).Value = 103.7191927446945
$ws.Cells.Item(16, 7).Value = 0.1290471052555283
$ws.Cells.Item(16, 8).Value = 2.983777742016493
$ws.Cells.Item(17, 7).Value = 0.1264717999989978
$ws.Cells.Item(17, 8).Value = -9.824344333390025
$ws.Cells.Item(18, 7).Value = 0.1109481945788235
$ws.Cells.Item(18, 8).Value = -11.05125341928895
$ws.Cells.Item(19, 7).Value = 0.1423782222973404
$ws.Cells.Item(19, 8).Value = 6.879449778859779
$ws.Cells.Item(20, 7).Value = 0.05238202104627122
$ws.Cells.Item(20, 8).Value = 52.55664344221179
$ws.Cells.Item(21, 7).Value = 0.06256955005424401
$ws.Cells.Item(21, 8).Value = 7.804442044647271
$ws.Cells.Item(22, 7).Value = -0.0630859155962649
$ws.Cells.Item(22, 8).Value = 20.9924755182562
$ws.Cells.Item(23, 7).Value = -0.07537832687715189
$ws.Cells.Item(23, 8).Value = -20.51578572357087
$ws.Cells.Item(24, 7).Value = 0.1214710005323355
$ws.Cells.Item(24, 8).Value = 2.841725791061934
$ws.Cells.Item(25, 7).Value = 0.1422446092902181
$ws.Cells.Item(25, 8).Value = 12.74008511208895
$ws.Cells.Item(26, 7).Value = 0.06864390289273559
$ws.Cells.Item(26, 8).Value = 38.10508904423892
$ws.Cells.Item(27, 7).Value = 0.07550719331837963
$ws.Cells.Item(27, 8).Value = -12.88627012599463
$ws.Cells.Item(28, 7).Value = -0.04274480759868867
$ws.Cells.Item(28, 8).Value = 32.78488053226422
$ws.Cells.Item(29, 7).Value = -0.07919836721646518
$ws.Cells.Item(29, 8).Value = -11.27707687750131
$ws.Cells.Item(30, 7).Value = 0.0519509502697816
$ws.Cells.Item(30, 8).Value = -18.45475156148083
$ws.Cells.Item(31, 7).Value = 0.04956714101311402
$ws.Cells.Item(31, 8).Value = -18.17965730470301
$ws.Cells.Item(32, 7).Value = 0.0839831733104256
$ws.Cells.Item(32, 8).Value = -14.53310798960525
$ws.Cells.Item(33, 7).Value = 0.1049957188706221
$ws.Cells.Item(33, 8).Value = 27.60052879821003
$ws.Cells.Item(34, 7).Value = 0.003522969919607256
$ws.Cells.Item(34, 8).Value = -86.47883047946455
$ws.Cells.Item(35, 7).Value = 0.01655421379083321
$ws.Cells.Item(35, 8).Value = 247.6721944099407
$ws.Cells.Item(36, 7).Value = 0.004052127665321132
$ws.Cells.Item(36, 8).Value = 652.8937848241521
$ws.Cells.Item(37, 7).Value = 0.01506212807822087
$ws.Cells.Item(37, 8).Value = 219.976184576555
$ws.Cells.Item(38, 7).Value = 0.1024096917410322
$ws.Cells.Item(38, 8).Value = -4.519577840260474
$ws.Cells.Item(39, 7).Value = 0.08510408494815715
$ws.Cells.Item(39, 8).Value = -0.6522050851806585
$ws.Cells.Item(40, 7).Value = 0.02740406013038355
$ws.Cells.Item(40, 8).Value = 822.6157021595059
$ws.Cells.Item(41, 7).Value = 0.04863942471396408
$ws.Cells.Item(41, 8).Value = 224.3331807454239
$ws.Cells.Item(42, 7).Value = 0.1093641216223649
$ws.Cells.Item(42, 8).Value = 8.350904655399097
$ws.Cells.Item(43, 7).Value = 0.1219161033385663
$ws.Cells.Item(43, 8).Value = 1.474557780853309
$ws.Cells.Item(44, 7).Value = 0.01901941646131767
$ws.Cells.Item(44, 8).Value = -46.70579318933977
$ws.Cells.Item(45, 7).Value = 0.02082640547857481
$ws.Cells.Item(45, 8).Value = 27.22164877817643
$ws.Cells.Item(46, 7).Value = 0.0558076554697642
$ws.Cells.Item(46, 8).Value = 54.00092716026504
$ws.Cells.Item(47, 7).Value = 0.05007300035638575
$ws.Cells.Item(47, 8).Value = -0.7282141641283327
$ws.Cells.Item(48, 7).Value = 0.05937936393623345
$ws.Cells.Item(48, 8).Value = 38.79746723028645
$ws.Cells.Item(49, 7).Value = 0.05582788613737956
$ws.Cells.Item(49, 8).Value = -19.64447400153448
$ws.Cells.Item(50, 7).Value = 0.01047673758537307
$ws.Cells.Item(50, 8).Value = -39.34513046419259
$ws.Cells.Item(51, 7).Value = 0.02878588105931734
$ws.Cells.Item(51, 8).Value = 47.84955457558868
$ws.Cells.Item(52, 7).Value = -0.1005095727870571
$ws.Cells.Item(52, 8).Value = 2.90849051334353
$ws.Cells.Item(53, 7).Value = -0.09712611110614215
$ws.Cells.Item(53, 8).Value = -5.166046699938506
$ws.Cells.Item(54, 7).Value = 0.0811168749654911
$ws.Cells.Item(54, 8).Value = 10.93260805072724
$ws.Cells.Item(55, 7).Value = 0.0998975973452785
$ws.Cells.Item(55, 8).Value = 61.25088761801489
$ws.Cells.Item(56, 7).Value = 0.04708117306956108
$ws.Cells.Item(56, 8).Value = 34.55910101441503
$ws.Cells.Item(57, 7).Value = 0.01869128610458377
$ws.Cells.Item(57, 8).Value = 223.7417288684268
$ws.Cells.Item(58, 7).Value = 0.05634468125061601
$ws.Cells.Item(58, 8).Value = 125.2831423260615
$ws.Cells.Item(59, 7).Value = 0.03190254778729745
$ws.Cells.Item(59, 8).Value = 34.73141593213877
$ws.Cells.Item(60, 7).Value = 0.01505244141420111
$ws.Cells.Item(60, 8).Value = -53.60293408956986
$ws.Cells.Item(61, 7).Value = 0.03766304504648783
$ws.Cells.Item(61, 8).Value = 197.5435722016105
$ws.Cells.Item(62, 7).Value = 0.05507919843926328
$ws.Cells.Item(62, 8).Value = -8.752800415095932
$ws.Cells.Item(63, 7).Value = 0.04928731244212137
$ws.Cells.Item(63, 8).Value = 51.23681709549148
$ws.Cells.Item(64, 7).Value = 0.03322875230072785
$ws.Cells.Item(64, 8).Value = -18.00670215516127
$ws.Cells.Item(65, 7).Value = 0.0615934778436878
$ws.Cells.Item(65, 8).Value = 9.866869670941997
$ws.Cells.Item(66, 7).Value = 0.08100674193448053
$ws.Cells.Item(66, 8).Value = -13.41187830260216
$ws.Cells.Item(67, 7).Value = 0.1050022977493404
$ws.Cells.Item(67, 8).Value = -9.046913095196048
$ws.Cells.Item(68, 7).Value = -0.02116900701428254
$ws.Cells.Item(68, 8).Value = 39.2576364184059
$ws.Cells.Item(69, 7).Value = -0.01339087482422184
$ws.Cells.Item(69, 8).Value = 36.90053875747415
$ws.Cells.Item(70, 7).Value = 0.08851271618061617
$ws.Cells.Item(70, 8).Value = -4.45132593605931
$ws.Cells.Item(71, 7).Value = 0.06097822295221921
$ws.Cells.Item(71, 8).Value = -33.14436169828681
$ws.Cells.Item(72, 7).Value = -0.05998773815280403
$ws.Cells.Item(72, 8).Value = -6.966192098130263
$ws.Cells.Item(73, 7).Value = -0.07827481482809583
$ws.Cells.Item(73, 8).Value = -6.117226853106361
$ws.Cells.Item(74, 7).Value = 0.1275261980563029
$ws.Cells.Item(74, 8).Value = 27.59406943110226
$ws.Cells.Item(75, 7).Value = 0.1186797855745743
$ws.Cells.Item(75, 8).Value = 21.84279693717349
$ws.Cells.Item(76, 7).Value = 0.02512954906205114
$ws.Cells.Item(76, 8).Value = -1.724829801273578
$ws.Cells.Item(77, 7).Value = 0.02164218710591856
$ws.Cells.Item(77, 8).Value = 53.39234192094349
$ws.Cells.Item(78, 7).Value = 0.07738345381488025
$ws.Cells.Item(78, 8).Value = 20.39063988223576
$ws.Cells.Item(79, 7).Value = 0.1069891698533084
$ws.Cells.Item(79, 8).Value = 39.46517527856352
$ws.Cells.Item(80, 7).Value = -0.1446587338575151
$ws.Cells.Item(80, 8).Value = 12.65012413601456
$ws.Cells.Item(81, 7).Value = -0.1893285503452016
$ws.Cells.Item(81, 8).Value = 9.883033496360735
$ws.Cells.Item(82, 7).Value = 0.1316083541619547
$ws.Cells.Item(82, 8).Value = 14.74604627285313
$ws.Cells.Item(83, 7).Value = 0.1605878061021735
$ws.Cells.Item(83, 8).Value = -9.773115054354134
$ws.Cells.Item(84, 7).Value = 0.06832131258123465
$ws.Cells.Item(84, 8).Value = 186.610950830866
$ws.Cells.Item(85, 7).Value = 0.07250504931734673
$ws.Cells.Item(85, 8).Value = 17.74899266557595
